$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 18 so everything below shifts up by one row.
$ws.Rows.Item(18).Delete()
